$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3:E11").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "(0.5989,0.30008)"
$ws.Range("E2").Value = "-0.0163"

# Row 3
$ws.Range("B3").Value = "(0.3082,0.27783)"
$ws.Range("C3").Value = "(0.30913,0.27897)"
$ws.Range("D3").Value = "-0.11603"
$ws.Range("E3").Value = "-0.2278"

# Row 4
$ws.Range("B4").Value = "(0.05246,0.01999)"
$ws.Range("C4").Value = "(0.05195,0.01774)"
$ws.Range("D4").Value = "0.06305"
$ws.Range("E4").Value = "0.45022"

# Row 5
$ws.Range("B5").Value = "(0.39949,0.19011)"
$ws.Range("C5").Value = "(0.39838,0.1891)"
$ws.Range("D5").Value = "0.13959"
$ws.Range("E5").Value = "0.20268"

# Row 6
$ws.Range("B6").Value = "(0.78175,0.31904)"
$ws.Range("C6").Value = "(0.78045,0.31426)"
$ws.Range("D6").Value = "0.16254"
$ws.Range("E6").Value = "0.95531"

# Row 7
$ws.Range("B7").Value = "(0.2981,0.47984)"
$ws.Range("C7").Value = "(0.2977,0.47797)"
$ws.Range("D7").Value = "0.04988"
$ws.Range("E7").Value = "0.37359"

# Row 8
$ws.Range("B8").Value = "(0.50199,0.23276)"
$ws.Range("C8").Value = "(0.5032,0.23243)"
$ws.Range("D8").Value = "-0.15121"
$ws.Range("E8").Value = "0.06633"

# Row 9
$ws.Range("B9").Value = "(0.52727,0.4496)"
$ws.Range("C9").Value = "(0.52899,0.44822)"
$ws.Range("D9").Value = "-0.21424"
$ws.Range("E9").Value = "0.27622"

# Row 10
$ws.Range("B10").Value = "(0.72516,0.45098)"
$ws.Range("C10").Value = "(0.72779,0.45464)"
$ws.Range("D10").Value = "-0.32872"
$ws.Range("E10").Value = "-0.73229"

# Row 11
$ws.Range("B11").Value = "(0.56015,0.14296)"
$ws.Range("C11").Value = "(0.56019,0.14323)"
$ws.Range("D11").Value = "-0.0053"
$ws.Range("E11").Value = "-0.05451"
